$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 66679.336
$ws.Range("I21").Value = 66679.336
$ws.Range("K21").Value = 66679.336
$ws.Range("M21").Value = -66211.336

$ws.Range("H23").Value = 66679.336
$ws.Range("I23").Value = 66679.336
$ws.Range("K23").Value = 66679.336
$ws.Range("M23").Value = -66445.336

$ws.Range("H40").Value = 3831.3333
$ws.Range("I40").Value = 3322.6667
$ws.Range("J40").Value = 4340
$ws.Range("K40").Value = 3322.6667
$ws.Range("L40").Value = 4340
$ws.Range("M40").Value = -3147.6667
$ws.Range("N40").Value = -4690

$ws.Range("H53").Value = 471.3125
$ws.Range("I53").Value = 136.85715
$ws.Range("K53").Value = 136.85715
$ws.Range("M53").Value = 500.14285

$ws.Range("H55").Value = 2017.2727
$ws.Range("J55").Value = 6899
$ws.Range("L55").Value = 6899
$ws.Range("N55").Value = -7327

$ws.Range("H64").Value = 4794.737
$ws.Range("J64").Value = 4794.737
$ws.Range("L64").Value = 4794.737
$ws.Range("N64").Value = -5290.737

$ws.Range("H67").Value = 4794.737
$ws.Range("J67").Value = 4794.737
$ws.Range("L67").Value = 4794.737
$ws.Range("N67").Value = -6510.737

$ws.Range("H107").Value = 3252.4546
$ws.Range("I107").Value = 2611.1428
$ws.Range("J107").Value = 4374.75
$ws.Range("K107").Value = 2611.1428
$ws.Range("L107").Value = 4374.75
$ws.Range("M107").Value = -691.1428000000001
$ws.Range("N107").Value = -8214.75

$ws.Range("H109").Value = 62495
$ws.Range("J109").Value = 62495
$ws.Range("L109").Value = 62495
$ws.Range("N109").Value = -65269

$ws.Range("H112").Value = 2562.5334
$ws.Range("J112").Value = 2495.2307
$ws.Range("L112").Value = 7485.6921
$ws.Range("N112").Value = -9701.6921

$ws.Range("H129").Value = 1535
$ws.Range("I129").Value = 1342
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 4026
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 974
$ws.Range("N129").Value = -17500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17866530
$ws.Range("I32").Value = 18528186
$ws.Range("K32").Value = 18528186
$ws.Range("M32").Value = -18527899

$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H74").Value = 17872564
$ws.Range("I74").Value = 31252676
$ws.Range("J74").Value = 32414.166
$ws.Range("K74").Value = 31252676
$ws.Range("L74").Value = 32414.166
$ws.Range("M74").Value = -31251802
$ws.Range("N74").Value = -34162.166

$ws.Range("H77").Value = 17872564
$ws.Range("I77").Value = 31252676
$ws.Range("J77").Value = 32414.166
$ws.Range("K77").Value = 156263380
$ws.Range("L77").Value = 162070.83
$ws.Range("M77").Value = -156259012
$ws.Range("N77").Value = -170806.83

$ws.Range("H102").Value = 2151.95
$ws.Range("I102").Value = 1891.1111
$ws.Range("K102").Value = 1891.1111
$ws.Range("M102").Value = -269.1111000000001

$ws.Range("H122").Value = 1934.1666
$ws.Range("I122").Value = 1550
$ws.Range("K122").Value = 4650
$ws.Range("M122").Value = -2200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2298.1
$ws.Range("J105").Value = 2896.6667
$ws.Range("L105").Value = 2896.6667
$ws.Range("N105").Value = -6390.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 299
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 839742.8
$ws.Range("I31").Value = 12074.786
$ws.Range("J31").Value = 1667410.8
$ws.Range("K31").Value = 12074.786
$ws.Range("L31").Value = 1667410.8
$ws.Range("M31").Value = -11779.786
$ws.Range("N31").Value = -1668000.8

$ws.Range("H34").Value = 839742.8
$ws.Range("I34").Value = 12074.786
$ws.Range("J34").Value = 1667410.8
$ws.Range("K34").Value = 12074.786
$ws.Range("L34").Value = 1667410.8
$ws.Range("M34").Value = -11872.786
$ws.Range("N34").Value = -1667814.8

$ws.Range("H62").Value = 2982.5
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876

$ws.Range("H65").Value = 2982.5
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380

$ws.Range("H107").Value = 741.84
$ws.Range("I107").Value = 558.25
$ws.Range("J107").Value = 1476.2
$ws.Range("K107").Value = 558.25
$ws.Range("L107").Value = 1476.2
$ws.Range("M107").Value = 1361.75
$ws.Range("N107").Value = -5316.2

$ws.Range("H122").Value = 4213.75
$ws.Range("I122").Value = 2242.1
$ws.Range("J122").Value = 7499.8335
$ws.Range("K122").Value = 6726.299999999999
$ws.Range("L122").Value = 22499.5005
$ws.Range("M122").Value = -4276.299999999999
$ws.Range("N122").Value = -27399.5005

$ws.Range("H134").Value = 437868.44
$ws.Range("I134").Value = 556720.75
$ws.Range("J134").Value = 10000.2
$ws.Range("K134").Value = 1670162.25
$ws.Range("L134").Value = 30000.6
$ws.Range("M134").Value = -1667627.25
$ws.Range("N134").Value = -35070.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11836153
$ws.Range("J4").Value = 25262516
$ws.Range("L4").Value = 75787548
$ws.Range("N4").Value = -75787772

$ws.Range("H34").Value = 3493.375
$ws.Range("J34").Value = 6249.75
$ws.Range("L34").Value = 18749.25
$ws.Range("N34").Value = -18917.25

$ws.Range("H48").Value = 1394
$ws.Range("J48").Value = 1394
$ws.Range("L48").Value = 4182
$ws.Range("N48").Value = -4682

$ws.Range("H55").Value = 12150
$ws.Range("J55").Value = 12150
$ws.Range("L55").Value = 36450
$ws.Range("N55").Value = -36804

$ws.Range("H92").Value = 1112679.9
$ws.Range("I92").Value = 1668237.1
$ws.Range("K92").Value = 5004711.300000001
$ws.Range("M92").Value = -5003463.300000001

$ws.Range("H117").Value = 725.6667
$ws.Range("J117").Value = 963.5
$ws.Range("L117").Value = 2890.5
$ws.Range("N117").Value = -9774.5

$ws.Range("H126").Value = 6000

$ws.Range("H129").Value = 30394424
$ws.Range("I129").Value = 1072
$ws.Range("J129").Value = 55722216
$ws.Range("K129").Value = 3216
$ws.Range("L129").Value = 167166648
$ws.Range("M129").Value = 1784
$ws.Range("N129").Value = -167176648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2538.6
$ws.Range("I102").Value = 2542.889
$ws.Range("K102").Value = 2542.889
$ws.Range("M102").Value = -920.8890000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2799.8
$ws.Range("I22").Value = 3250
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 3250
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -2955
$ws.Range("N22").Value = -1589

$ws.Range("H27").Value = 2799.8
$ws.Range("I27").Value = 3250
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 3250
$ws.Range("L27").Value = 999
$ws.Range("M27").Value = -3143
$ws.Range("N27").Value = -1213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2233.7058
$ws.Range("I96").Value = 1647.9
$ws.Range("J96").Value = 3070.5715
$ws.Range("K96").Value = 1647.9
$ws.Range("L96").Value = 3070.5715
$ws.Range("M96").Value = -274.9000000000001
$ws.Range("N96").Value = -5816.5715

$ws.Range("H100").Value = 1944.0408
$ws.Range("I100").Value = 1894.6222
$ws.Range("K100").Value = 3789.2444
$ws.Range("M100").Value = -3248.2444

$ws.Range("H107").Value = 41668480
$ws.Range("I107").Value = 62502220
$ws.Range("K107").Value = 187506660
$ws.Range("M107").Value = -187504740

$ws.Range("H122").Value = 3792.5
$ws.Range("I122").Value = 2620.8
$ws.Range("J122").Value = 6455.4546
$ws.Range("K122").Value = 7862.400000000001
$ws.Range("L122").Value = 19366.3638
$ws.Range("M122").Value = -5412.400000000001
$ws.Range("N122").Value = -24266.3638
